# Remove the standalone "Miquéias" (italic) heading paragraph that sits
# directly under the "MIC" Heading2 paragraph. The "MIC" paragraph itself
# is left untouched; only the following paragraph (which contains nothing
# but the italic book-name run) is deleted in its entirety, including its
# paragraph mark, so the paragraph that used to follow it moves up.

$d = $word.ActiveDocument

# Locate the italic "Miquéias" sub-heading that immediately follows "MIC".
# NOTE: $d.Content returns a fresh whole-document Range every time it is
# evaluated, so it must be captured in a variable before calling Find -
# Find.Execute mutates that specific Range object (collapsing it to the
# match) rather than any document-level state.
$rng = $d.Content
$found = $rng.Find.Execute("Miquéias", $true, $false, $false, $false,
                            $false, $true, 1, $false, "", 0)

if ($found -and $rng.Italic) {
    $hitStart = $rng.Start

    # Resolve which paragraph contains the found text.
    $targetIndex = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $candidate = $d.Paragraphs.Item($i)
        if ($candidate.Range.Start -le $hitStart -and $hitStart -lt $candidate.Range.End) {
            $targetIndex = $i
        }
    }

    if ($targetIndex -gt 0) {
        $targetPara = $d.Paragraphs.Item($targetIndex)
        $afterPara = $d.Paragraphs.Item($targetIndex + 1)

        # Delete from the start of the target paragraph through the start of
        # the next paragraph so the paragraph mark itself is removed too
        # (otherwise an empty paragraph would be left behind).
        $killRange = $d.Range($targetPara.Range.Start, $afterPara.Range.Start)
        $killRange.Delete()
    }
}
